$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 1.0.1 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text fix (applies to all 4 test-case blocks, shared string)
$newPrecondition = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B8").Value = $newPrecondition
$ws.Range("B17").Value = $newPrecondition
$ws.Range("B25").Value = $newPrecondition
$ws.Range("B33").Value = $newPrecondition

# TC1 step 2 expected results: accent corrections
$ws.Range("D11").Value = "SYSTEM Exibe a lista de diárias (solicitações) aptas para pagamento ordenado pelo número da diária em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de liquidação (após registrar o empenho)."

# TC2 step 2 (row 20): swap content with what used to be TC3 step 2
$ws.Range("B20").Value = "Chefe Clica para realizar a liquidação."
$ws.Range("D20").Value = "SYSTEM Apresenta a tela de Registrar Liquidações."

# TC3 step 2 (row 28): swap content with what used to be TC2 step 2 (with wording tweak)
$ws.Range("B28").Value = "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D28").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# TC4 step 2 (row 36): add trailing period
$ws.Range("D36").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."
